# Atualização de bases das ligas, do dia: 03-03-2024 às 00:35
# Fills in match results (FTHG, FTAG, FTR) and closing odds / P&L columns
# for rows 121 and 122 (matches that have now been played), and updates
# the closing odds for row 123.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Row 121 (Hajduk Split vs HNK Gorica) ----
$ws.Range("H121").Value = 2
$ws.Range("I121").Value = 3
$ws.Range("J121").Value = "A"

$ws.Range("N121").Value = 2.625
$ws.Range("O121").Value = 2.9
$ws.Range("P121").Value = 2.75
$ws.Range("R121").Value = 1.875
$ws.Range("S121").Value = 1.975
$ws.Range("U121").Value = 2.1
$ws.Range("V121").Value = 1.775
$ws.Range("W121").Value = -1
$ws.Range("X121").Value = -1
$ws.Range("Y121").Value = 1.75
$ws.Range("Z121").Value = -1
$ws.Range("AA121").Value = 0.9750000000000001
$ws.Range("AB121").Value = 1.1
$ws.Range("AC121").Value = -1

# ---- Row 122 (Istra 1961 vs NK Rudes) ----
$ws.Range("H122").Value = 1
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = "H"

$ws.Range("N122").Value = 1.285
$ws.Range("P122").Value = 10
$ws.Range("R122").Value = 2
$ws.Range("S122").Value = 1.85
$ws.Range("U122").Value = 1.975
$ws.Range("V122").Value = 1.875
$ws.Range("W122").Value = 0.2849999999999999
$ws.Range("X122").Value = -1
$ws.Range("Y122").Value = -1
$ws.Range("Z122").Value = -1
$ws.Range("AA122").Value = 0.8500000000000001
$ws.Range("AB122").Value = -1
$ws.Range("AC122").Value = 0.875

# ---- Row 123 (closing odds update only) ----
$ws.Range("N123").Value = 1.285
$ws.Range("T123").Value = 2.75
$ws.Range("U123").Value = 2.05
$ws.Range("V123").Value = 1.8

$wb.Save()
